$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark; Word will leave it at the very
# end of the document's last edit once we're done (right after the
# "Worked?" entry typed into the final table row).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count

# Fill in "Y" for every data row's "Worked?" column (column 4), skipping
# the header row (row 1).
for ($i = 2; $i -le $rowCount; $i++) {
    $cell = $t.Cell($i, 4)
    $cell.Range.Text = "Y"
    $cell.Range.Font.Size = 12
}

# Re-create the _GoBack bookmark collapsed right after the "Y" just
# typed into the last row's "Worked?" cell. Word records _GoBack at the
# site of the most recent edit, which is this final cell.
#
# A collapsed range sitting exactly on a paragraph-mark boundary can't
# be fed straight into Bookmarks.Add, so temporarily insert a marker
# character after the "Y", bookmark across it (a real, non-boundary
# range), then delete the marker - the bookmark collapses in place and
# survives.
$lastCell = $t.Cell($rowCount, 4)
$insertionPoint = $lastCell.Range.Duplicate
$insertionPoint.Collapse(0) | Out-Null
$insertionPoint.MoveEnd(1, -1) | Out-Null
$insertionPoint.Collapse(0) | Out-Null
$insertionPoint.InsertAfter("~") | Out-Null
$markerStart = $insertionPoint.Start
$markerRange = $d.Range($markerStart, $markerStart + 1)
$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null
$markerRange2 = $d.Range($markerStart, $markerStart + 1)
$markerRange2.Text = ""
